$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D6 target cluster label to "Resolving-Mac" before removing the Neutrophils row,
# so the surviving shared-string slot repoints instead of leaving a stale one behind.
$ws.Range("D6").Value = "Resolving-Mac"

# Refresh the recomputed TPM-based metrics for the remaining target clusters.
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08785
$ws.Range("H2").Value = 0.26355
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.70092
$ws.Range("N2").Value = 5.10276
$ws.Range("O2").Value = 0.06397439760344623
$ws.Range("P2").Value = 0.06397439760344621
$ws.Range("Q2").Value = 0.149425822
$ws.Range("R2").Value = 1.344832398
$ws.Range("S2").Value = 0.06397439760344623
$ws.Range("T2").Value = 0.06397439760344621
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08785
$ws.Range("H3").Value = 0.26355
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.74650133333333
$ws.Range("N3").Value = 35.239504
$ws.Range("O3").Value = 0.4418052270230686
$ws.Range("P3").Value = 0.4418052270230685
$ws.Range("Q3").Value = 1.031930142133333
$ws.Range("R3").Value = 9.287371279199998
$ws.Range("S3").Value = 0.4418052270230686
$ws.Range("T3").Value = 0.4418052270230685
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08785
$ws.Range("H4").Value = 0.26355
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.040268
$ws.Range("N4").Value = 0.120804
$ws.Range("O4").Value = 0.001514545682745557
$ws.Range("P4").Value = 0.001514545682745557
$ws.Range("Q4").Value = 0.0035375438
$ws.Range("R4").Value = 0.0318378942
$ws.Range("S4").Value = 0.001514545682745557
$ws.Range("T4").Value = 0.001514545682745557
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08785
$ws.Range("H5").Value = 0.26355
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.96401166666667
$ws.Range("N5").Value = 38.892035
$ws.Range("O5").Value = 0.4875977923118364
$ws.Range("P5").Value = 0.4875977923118364
$ws.Range("Q5").Value = 1.138888424916667
$ws.Range("R5").Value = 10.24999582425
$ws.Range("S5").Value = 0.4875977923118364
$ws.Range("T5").Value = 0.4875977923118364
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08785
$ws.Range("H6").Value = 0.26355
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.13581
$ws.Range("N6").Value = 0.40743
$ws.Range("O6").Value = 0.005108037378903201
$ws.Range("P6").Value = 0.0051080373789032
$ws.Range("Q6").Value = 0.0119309085
$ws.Range("R6").Value = 0.1073781765
$ws.Range("S6").Value = 0.005108037378903201
$ws.Range("T6").Value = 0.0051080373789032

# Drop the old Neutrophils row (row 7) entirely - no longer present in the new data.
$ws.Rows("7:7").Delete()
